$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-CellText($row, $col, $old, $new) {
    $found = $t.Cell($row, $col).Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: replacement failed for row=$row col=$col old=$old"
    }
}

# Row 1 (first data row of table)
Replace-CellText 1 1 "70÷9=" "33÷8="
Replace-CellText 1 2 "39÷9=" "31÷6="
Replace-CellText 1 3 "99÷5=" "94÷7="
Replace-CellText 1 4 "22÷5=" "41÷8="
Replace-CellText 1 5 "23÷2=" "22÷2="

# Row 5 (second data row of table)
Replace-CellText 5 1 "91÷4=" "94÷6="
Replace-CellText 5 2 "55÷5=" "97÷2="
Replace-CellText 5 3 "79÷9=" "33÷6="
Replace-CellText 5 4 "32÷6=" "65÷9="
Replace-CellText 5 5 "62÷7=" "60÷3="

# Row 9 (third data row of table)
Replace-CellText 9 1 "63÷9=" "58÷6="
Replace-CellText 9 2 "98÷6=" "86÷3="
Replace-CellText 9 3 "23÷9=" "53÷9="
Replace-CellText 9 4 "58÷9=" "95÷3="
Replace-CellText 9 5 "21÷9=" "67÷2="

# Row 13 (fourth data row of table)
Replace-CellText 13 1 "56÷3=" "63÷9="
Replace-CellText 13 2 "17÷7=" "98÷9="
Replace-CellText 13 3 "13÷5=" "81÷9="
Replace-CellText 13 4 "47÷3=" "58÷3="
Replace-CellText 13 5 "67÷4=" "23÷9="

# Row 17 (fifth data row of table)
Replace-CellText 17 1 "30÷3=" "93÷4="
Replace-CellText 17 2 "36÷6=" "53÷2="
Replace-CellText 17 3 "53÷6=" "28÷5="
Replace-CellText 17 4 "64÷7=" "74÷6="
Replace-CellText 17 5 "28÷9=" "72÷8="
